# Atualização RDD 10 - Liga Eliminação
# Ajuste e atualização da pontuação da Rodada 10 na Liga Eliminação.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Goleiro ---
$ws.Range("D2").Value = 48    # Anthoni       47 -> 48
$ws.Range("D3").Value = 46    # Weverton      44 -> 46
# Ronaldo (row 4) unchanged

# --- Zagueiro ---
$ws.Range("D5").Value = 135   # Léo Ortiz     134 -> 135

# Ignácio dropped out, replaced at the top of the Zagueiro group by
# Pedro Henrique (RBB); remaining Zagueiros shift down one spot and
# Vitão drops out of the ranking entirely.
$ws.Range("A6").Value = "Pedro Henrique"
$ws.Range("C6").Value = "RBB"
$ws.Range("D6").Value = 42

$ws.Range("A7").Value = "Ignácio"
$ws.Range("C7").Value = "FLU"
$ws.Range("D7").Value = 32

$ws.Range("A8").Value = "Junior Alonso"
$ws.Range("C8").Value = "CAM"
$ws.Range("D8").Value = 31

# --- Lateral ---
$ws.Range("D9").Value = 118   # Juninho Capixaba  88 -> 118
$ws.Range("D10").Value = 79   # Bernabéi          78 -> 79
# Escobar (row 11) and Wesley (row 12) unchanged

# --- Meia ---
$ws.Range("D13").Value = 157  # Arias         129 -> 157
$ws.Range("D14").Value = 131  # Arrascaeta    127 -> 131
$ws.Range("D15").Value = 99   # Alan Patrick  98 -> 99

# Jhon Jhon moves above Rubens in the Meia group
$ws.Range("A16").Value = "Jhon Jhon"
$ws.Range("C16").Value = "RBB"
$ws.Range("D16").Value = 58

$ws.Range("A17").Value = "Rubens"
$ws.Range("C17").Value = "CAM"
$ws.Range("D17").Value = 52

# --- Atacante ---
$ws.Range("D18").Value = 114  # Yuri Alberto  113 -> 114
$ws.Range("D19").Value = 68   # Hulk          62 -> 68
# Igor Jesus (row 20) unchanged

# Wesley drops out of the Atacante ranking, replaced by Eduardo Sasha
$ws.Range("A21").Value = "Eduardo Sasha"
$ws.Range("C21").Value = "RBB"
$ws.Range("D21").Value = 53

$ws.Range("D22").Value = 53   # Estêvão       51 -> 53

$wb.Save()
